# Update the D12:H12 F1-score row ("Metabuli") on each results sheet with
# the new Metabuli version results, as described in the commit.
# Sheet5 (16S / Lutjanidae) is left unchanged.

$wb = $excel.ActiveWorkbook

function Set-MetabuliRow($SheetIndex, $D, $E, $F, $G, $H) {
    $ws = $wb.Worksheets.Item($SheetIndex)
    $ws.Range("D12").Value = $D
    $ws.Range("E12").Value = $E
    $ws.Range("F12").Value = $F
    $ws.Range("G12").Value = $G
    $ws.Range("H12").Value = $H
}

# Sheet1 - 12S, 100 Australian species
Set-MetabuliRow 1 0.8115942028985508 0.6511627906976745 0.7225806451612903 0.7734806629834254 0.5656565656565656

# Sheet2 - 16S, 100 Australian species
Set-MetabuliRow 2 0.7464788732394366 0.654320987654321 0.6973684210526315 0.726027397260274 0.5353535353535354

# Sheet3 - COI, 100 Australian species
Set-MetabuliRow 3 0.9491525423728814 0.8235294117647058 0.8818897637795277 0.9210526315789472 0.8484848484848485

# Sheet4 - 12S, Lutjanidae
Set-MetabuliRow 4 0.9411764705882353 0.7619047619047619 0.8421052631578947 0.898876404494382 0.75

# Sheet5 - 16S, Lutjanidae -> unchanged, no edits

# Sheet6 - COI, Lutjanidae
Set-MetabuliRow 6 1 0.92 0.9583333333333334 0.9829059829059831 0.9259259259259259

# Sheet7 - 12S, Wadjemup (H12 unchanged at 0.6470588235294118)
Set-MetabuliRow 7 0.8552631578947368 0.7222222222222222 0.783132530120482 0.8248730964467005 0.6470588235294118

# Sheet8 - 16S, Wadjemup
Set-MetabuliRow 8 0.7931034482758621 0.7666666666666667 0.7796610169491527 0.7876712328767123 0.6517857142857143

# Sheet9 - COI, Wadjemup
Set-MetabuliRow 9 0.9782608695652174 0.7964601769911505 0.8780487804878049 0.9355509355509356 0.7863247863247863
